$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J: existing J column shifts to K, existing K shifts to L
$ws.Columns("J:J").Insert()

# Row 1 (header) values are unaffected by the shift conceptually (still 0..9),
# restore J1/K1 and add the new L1 header value
$ws.Range("J1").Value = 8
$ws.Range("K1").Value = 9
$ws.Range("L1").Value = 10

# Fill the newly inserted J column (rows 2-12) with the new decimal values
$ws.Range("J2").Value = 0.219
$ws.Range("J3").Value = 0.181
$ws.Range("J4").Value = 0.148
$ws.Range("J5").Value = 0.133
$ws.Range("J6").Value = 0.108
$ws.Range("J7").Value = 0.078
$ws.Range("J8").Value = 0.059
$ws.Range("J9").Value = 0.043
$ws.Range("J10").Value = 0.042
$ws.Range("J11").Value = 0.042
$ws.Range("J12").Value = 0.041
